# Cotações atualizadas - 2025-11-30
# Append a new row (86) with the quotes for 2025-11-30 (serial date 45991),
# reusing the date style from the previous row (A85) so the new date cell
# keeps the same number format as the rest of column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (number format, etc.) of the last existing date cell (A85)
# onto the new date cell (A86) before writing its value.
$ws.Range("A85").Copy()
$ws.Range("A86").PasteSpecial(-4122)

$ws.Range("A86").Value = 45991
$ws.Range("B86").Value = "21,7883"
$ws.Range("C86").Value = "16,0515"
$ws.Range("D86").Value = "15,5122"
$ws.Range("E86").Value = "15,5122"
